$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits right
#    after the H1 title paragraph.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2) Replace the final paragraph ("Create a feature image...") with two new
#    paragraphs:
#      a) a bold paragraph containing the page title/meta "headline"
#      b) an italic paragraph containing the meta description text
#    (this mirrors the text that used to live at the top of the document).
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Delete()

$endOfDoc = $d.Paragraphs($d.Paragraphs.Count)
$insertPoint = $d.Range($endOfDoc.Range.End, $endOfDoc.Range.End)

$newParasXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Age of the Gods &#8211; Fate Sister for Free | Exciting Bonus Features</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Find out why Age of the Gods &#8211; Fate Sister is worth playing! Enjoy exciting bonus features and win progressive jackpots. Play for free!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($newParasXml) | Out-Null
